$wb = $excel.ActiveWorkbook

# --- Step 1: rename Sheet2 to May ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "May"

# --- Step 2: insert 2 new rows at the top of May's data table for May 30 and May 31 ---
$ws2.Rows("2:3").Insert()

# Copy date style (s="1") from the row below down into the new rows
$ws2.Range("A4").Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

$ws2.Range("A2").Value = 45443
$ws2.Range("B2").Value = 267
$ws2.Range("C2").Value = 24
$ws2.Range("D2").Value = 7
$ws2.Range("E2").Value = 236

$ws2.Range("A3").Value = 45442
$ws2.Range("B3").Value = 248
$ws2.Range("C3").Value = 28
$ws2.Range("D3").Value = 8
$ws2.Range("E3").Value = 212

# Fix the SUM formulas to cover the new rows (B2:B32 etc.)
$ws2.Range("B33:E33").Formula = "=SUM(B2:B32)"

Write-Host "May sheet updated"

# --- Step 3: add a new "Sheet3" sheet (June data) right after May ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item("Sheet1 (2)")
$ws3.Name = "Sheet3"
$ws3.Rows("2:32").Delete()

$ws3.Range("A2").Value = 45444
$ws3.Range("B2").Value = 209
$ws3.Range("C2").Value = 30
$ws3.Range("D2").Value = 6
$ws3.Range("E2").Value = 173

Write-Host "Sheet3 created"
